$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text) {
        $arr = $text -split ", "

        $hasSystem = $false
        for ($i = 0; $i -lt $arr.Length; $i++) {
            if ($arr[$i].Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem -and -not $arr[0].Equals("System")) {
            $newText = "System"
            for ($i = 0; $i -lt $arr.Length; $i++) {
                if (-not $arr[$i].Equals("System")) {
                    $newText = $newText + ", " + $arr[$i]
                }
            }
            $cell.Value = $newText
        }
    }
}
